$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target contents for columns B..F, rows 2-12, after the edit.
# $null means the cell must end up blank.
$data = @(
    @{ Row=2;  B="NSE:ANGELONE";   C="NSE:AGRITECH";   D=$null; E="NSE:MANAPPURAM"; F="NSE:GRANULES" }
    @{ Row=3;  B="NSE:GILLANDERS"; C="NSE:CORDSCABLE"; D=$null; E=$null;            F=$null }
    @{ Row=4;  B="NSE:GIPCL";      C="NSE:HIL";        D=$null; E=$null;            F=$null }
    @{ Row=5;  B="NSE:GRANULES";   C="NSE:HOMEFIRST";  D=$null; E=$null;            F=$null }
    @{ Row=6;  B="NSE:GROBTEA";    C="NSE:INDORAMA";   D=$null; E=$null;            F=$null }
    @{ Row=7;  B="NSE:KEYFINSERV"; C="NSE:KARURVYSYA"; D=$null; E=$null;            F=$null }
    @{ Row=8;  B="NSE:KITEX";      C="NSE:LICNFNHGP";  D=$null; E=$null;            F=$null }
    @{ Row=9;  B="NSE:PARAGMILK";  C="NSE:MAHLIFE";    D=$null; E=$null;            F=$null }
    @{ Row=10; B="NSE:PETRONET";   C="NSE:PALREDTEC";  D=$null; E=$null;            F=$null }
    @{ Row=11; B="NSE:SAFARI";     C="NSE:RITES";      D=$null; E=$null;            F=$null }
    @{ Row=12; B="NSE:SAHYADRI";   C=$null;            D=$null; E=$null;            F=$null }
)

foreach ($entry in $data) {
    $r = $entry.Row
    foreach ($colInfo in @(@{Col=2; Key="B"}, @{Col=3; Key="C"}, @{Col=4; Key="D"}, @{Col=5; Key="E"}, @{Col=6; Key="F"})) {
        $cell = $ws.Cells.Item($r, $colInfo.Col)
        $val = $entry[$colInfo.Key]
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value = $val
        }
    }
}

# Delete rows 13 through 29 entirely, shrinking the used range to A1:F12.
$ws.Range("A13:F29").EntireRow.Delete()
